$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "63.40", "5.50").
# Force Text format before assignment so Excel does not coerce the string
# into a Double and silently drop significant trailing zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.044.33"
$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.644.27"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("E5").Value = "  +0.90%  "

$ws.Range("E6").Value = "  +0.70%  "

$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("E8").Value = "  +0.53%  "

$ws.Range("E9").Value = "  +1.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  +0.29%  "

$ws.Range("E11").Value = "  +0.54%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.688.19"
$ws.Range("E12").Value = "  +3.24%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.873.89"
$ws.Range("E13").Value = "  +0.79%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.29"
$ws.Range("E14").Value = "  +1.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.545"
$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("E16").Value = "  +1.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.40"
$ws.Range("E17").Value = "  +1.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.060.64"
$ws.Range("E18").Value = "  +0.55%  "

$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "195.64"
$ws.Range("E20").Value = "  +1.31%  "

$ws.Range("E21").Value = "  -0.59%  "

$ws.Range("E23").Value = "  -0.23%  "

$ws.Range("E24").Value = "  +0.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.132"
$ws.Range("E25").Value = "  +5.01%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.88"
$ws.Range("E26").Value = "  +0.62%  "

$ws.Range("B27").Value = "BinanceUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.01"
$ws.Range("E27").Value = "  +0.58%  "

$ws.Range("E28").Value = "  +0.78%  "

$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.27"
$ws.Range("E32").Value = "  +1.50%  "

$ws.Range("E33").Value = "  -0.24%  "

$ws.Range("E34").Value = "  -2.42%  "

$ws.Range("E35").Value = "  +1.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.907"
$ws.Range("E36").Value = "  +0.43%  "

$ws.Range("B37").Value = "NeutrinoSystemBaseToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/fWTaNV_Ff+neutrinosystembasetoken-nsbt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "430.08"
$ws.Range("E37").Value = "  +20.14%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.134.28"
$ws.Range("E38").Value = "  -0.47%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.542"
$ws.Range("E39").Value = "  -1.39%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.47"
$ws.Range("E40").Value = "  -0.41%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0157"
$ws.Range("E41").Value = "  +0.20%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.50"
$ws.Range("E42").Value = "  +0.90%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.38"
$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.797"
$ws.Range("E44").Value = "  -0.94%  "

$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.782.71"
$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0117"
$ws.Range("E46").Value = "  +4.07%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.77"
$ws.Range("E47").Value = "  +0.82%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0526"
$ws.Range("E48").Value = "  +0.24%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.47"
$ws.Range("E49").Value = "  +0.35%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.79"
$ws.Range("E50").Value = "  +2.41%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.415"
$ws.Range("E51").Value = "  -0.03%  "
